# dados/ADD/Dados_ADD_PF/faturamento_diario.xlsx -- "atualizei dados bibi e add"
#
# Two changes to the daily-revenue table on Sheet1:
#   1) Two already-recorded daily totals ("dados bibi") were corrected:
#        - dia 1 / 07-2025 (row 2): total_venda 16332.98 -> 17296.13
#        - dia 4 / 07-2025 (row 5): total_venda 13554.45 -> 26369.89
#   2) A new daily entry ("add") for dia 7 / 07-2025 was added right after
#      the existing 07-2025 rows. Every row from the old row 6 through the
#      old last row (67) shifts down by one, so the sheet grows from
#      A1:E67 to A1:E68.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Corrected totals ---
$ws.Cells.Item(2, 2).Value = 17296.13
$ws.Cells.Item(5, 2).Value = 26369.89

# --- 2) New row for "dia 7" of 07/2025, inserted after the other 07/2025
#        rows; this pushes the rest of the table (old rows 6:67) down to
#        rows 7:68. ---
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).Value = 7
$ws.Cells.Item(6, 2).Value = 15800.46
$ws.Cells.Item(6, 3).Value = 7
$ws.Cells.Item(6, 4).Value = 2025
$ws.Cells.Item(6, 5).Value = "07/2025"
